$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 401
$ws.Range("E6").Value = 94
$ws.Range("G6").Value = 23.44139650872818
$ws.Range("H6").Value = 76.55860349127181
